# This script reproduces the edit described by the diff:
#  - Fills in two new "extraordinary" test cases in rows 11 and 12 of the
#    "Casos Extraordinarios" block (previously blank placeholder rows).
#  - Moves the existing "Ahorro programado estandar" block (the header row
#    and its three data rows) and the "Formulas utilizadas" block below it
#    down by two rows to make room (old rows 13-20 -> new rows 15-22),
#    leaving rows 13, 14 and 19 blank, and leaving the trailing row 23
#    (outside the edited area) untouched.
#  - Updates the active selection to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Capture the current contents (values/formulas) of the rows that need
#    to move, BEFORE we overwrite anything, so we can relocate them safely.
# ---------------------------------------------------------------------

# Row 13 (header row for the "Ahorro programado estandar" block)
$r13_A = $ws.Range("A13").Value()
$r13_B = $ws.Range("B13").Value()
$r13_C = $ws.Range("C13").Value()
$r13_D = $ws.Range("D13").Value()
$r13_E = $ws.Range("E13").Value()
$r13_F = $ws.Range("F13").Value()
$r13_G = $ws.Range("G13").Value()

# Row 14 (case 1 of that block)
$r14_A = $ws.Range("A14").Value()
$r14_B = $ws.Range("B14").Value()
$r14_C = $ws.Range("C14").Value()
$r14_D = $ws.Range("D14").Value()
$r14_E = $ws.Range("E14").Value()
$r14_F = $ws.Range("F14").Value()

# Row 15 (case 2 of that block)
$r15_A = $ws.Range("A15").Value()
$r15_B = $ws.Range("B15").Value()
$r15_C = $ws.Range("C15").Value()
$r15_D = $ws.Range("D15").Value()
$r15_E = $ws.Range("E15").Value()
$r15_F = $ws.Range("F15").Value()

# Row 16 (case 3 of that block)
$r16_A = $ws.Range("A16").Value()
$r16_B = $ws.Range("B16").Value()
$r16_C = $ws.Range("C16").Value()
$r16_D = $ws.Range("D16").Value()
$r16_E = $ws.Range("E16").Value()
$r16_F = $ws.Range("F16").Value()

# Row 18 ("Formulas utilizadas" title)
$r18_A = $ws.Range("A18").Value()

# Row 19 (first formula description line)
$r19_A = $ws.Range("A19").Value()
$r19_B = $ws.Range("B19").Value()

# Row 20 (second formula description line)
$r20_A = $ws.Range("A20").Value()
$r20_B = $ws.Range("B20").Value()

# ---------------------------------------------------------------------
# 2) Clear out the old block (formatting included) so we can rebuild it
#    cleanly two rows lower without leaving stray formatting behind.
# ---------------------------------------------------------------------
$ws.Range("A13:G20").Clear()

# ---------------------------------------------------------------------
# 3) Fill rows 11 and 12 with the two new extraordinary test cases.
# ---------------------------------------------------------------------

# Case 3: rate too low for 200 periods
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 200000
$ws.Range("C11").Value = 16000000
$ws.Range("D11").Value = 200
$ws.Range("E11").Value = 0.01
$ws.Range("F11").Value = "superior a 200 años"
$ws.Range("G11").Formula = "=IF(E11=0, B11 + C11*D11, B11*(1+E11)^D11 + C11*((1+E11)^D11-1)/E11)"

# Case 4: negative interest rate
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 205000
$ws.Range("C12").Value = 500000
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = -1
$ws.Range("F12").Value = "Tasa de interes en negativa"
$ws.Range("G12").Formula = "=IF(E12=0, B12 + C12*D12, B12*(1+E12)^D12 + C12*((1+E12)^D12-1)/E12)"

# ---------------------------------------------------------------------
# 4) Rows 13 and 14 stay blank (the gap created by the shift).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5) Rebuild the moved block two rows down (old 13 -> 15, .., old 16 -> 18).
# ---------------------------------------------------------------------

# Row 15: header row, every column bold (matching the original header's
# formatting, which overrides the currency column format in B/C/G).
$ws.Range("A15").Value = $r13_A
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").Copy()
$ws.Range("B15:G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B15").Value = $r13_B
$ws.Range("C15").Value = $r13_C
$ws.Range("D15").Value = $r13_D
$ws.Range("E15").Value = $r13_E
$ws.Range("F15").Value = $r13_F
$ws.Range("G15").Value = $r13_G

# Row 16: case 1 (bold label in column A only)
$ws.Range("A16").Value = $r14_A
$ws.Range("B16").Value = $r14_B
$ws.Range("C16").Value = $r14_C
$ws.Range("D16").Value = $r14_D
$ws.Range("E16").Value = $r14_E
$ws.Range("F16").Value = $r14_F
$ws.Range("G16").Formula = "=IF(E16=0, B16 + C16*D16, B16*(1+E16)^D16 + C16*((1+E16)^D16-1)/E16)"
$ws.Range("A16").Font.Bold = $true

# Row 17: case 2 (bold label in column A only)
$ws.Range("A17").Value = $r15_A
$ws.Range("B17").Value = $r15_B
$ws.Range("C17").Value = $r15_C
$ws.Range("D17").Value = $r15_D
$ws.Range("E17").Value = $r15_E
$ws.Range("F17").Value = $r15_F
$ws.Range("G17").Formula = "=IF(E17=0, B17 + C17*D17, B17*(1+E17)^D17 + C17*((1+E17)^D17-1)/E17)"
$ws.Range("A17").Font.Bold = $true

# Row 18: case 3 (bold label in column A only)
$ws.Range("A18").Value = $r16_A
$ws.Range("B18").Value = $r16_B
$ws.Range("C18").Value = $r16_C
$ws.Range("D18").Value = $r16_D
$ws.Range("E18").Value = $r16_E
$ws.Range("F18").Value = $r16_F
$ws.Range("G18").Formula = "=IF(E18=0, B18 + C18*D18, B18*(1+E18)^D18 + C18*((1+E18)^D18-1)/E18)"
$ws.Range("A18").Font.Bold = $true

# ---------------------------------------------------------------------
# 6) Row 19 stays blank (the gap); rebuild rows 20-22 from old rows 18-20.
# ---------------------------------------------------------------------

# Row 20: "Formulas utilizadas" title (bold)
$ws.Range("A20").Value = $r18_A
$ws.Range("A20").Font.Bold = $true

# Row 21 and 22: formula description lines (not bold)
$ws.Range("A21").Value = $r19_A
$ws.Range("B21").Value = $r19_B

$ws.Range("A22").Value = $r20_A
$ws.Range("B22").Value = $r20_B

# ---------------------------------------------------------------------
# 7) Update the active selection to F12, as in the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("F12").Select()
